$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep their original text formatting (Price/Volume columns
# contain text-like values such as "0.120" or "11.20" that Excel would otherwise
# auto-convert to numbers, stripping the formatting / trailing zeros).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.166.69"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.381.99"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.78%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.33"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.15"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.83%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.81%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +9.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.589"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "48.71"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +5.14%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +5.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "685.56"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.31%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.929.46"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.200.03"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.17%  "
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.120"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.88%  "
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.381.09"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.69"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.36"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.901"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.42"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.08"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.94%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "103.89"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +5.59%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.63%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.61"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.74%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +4.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.69"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.00%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.20"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.12%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +10.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "555.48"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.24%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.99"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.700.03"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.140"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +7.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.83"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.39%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.78%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.94%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.68"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.05%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0419"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.70%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.77%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.44%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.86%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +5.10%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.24%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.52%  "
